$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 132
$ws.Range("I2").Value = 341
$ws.Range("J2").Value = 1451
$ws.Range("K2").Value = 11
$ws.Range("L2").Value = 446
$ws.Range("M2").Value = 29
$ws.Range("N2").Value = 268
$ws.Range("P2").Value = 4
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 16
$ws.Range("S2").Value = 161
$ws.Range("T2").Value = 265
$ws.Range("U2").Value = 20
$ws.Range("V2").Value = 2273
$ws.Range("X2").Value = 2300
$ws.Range("Y2").Value = 10
$ws.Range("Z2").Value = 35
$ws.Range("AA2").Value = 14
